$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.636.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("E7").Value = "  -0.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3285"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7768"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07753"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.034"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("E17").Value = "  -1.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007982"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.656.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.083.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.642"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.549"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.201"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  +1.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.192"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("E31").Value = "  -1.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08763"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.79%  "

$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.871"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7152"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.197"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4888"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.083"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.736"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4203"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05923"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8848"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.20%  "
